$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.4437935948371887, -0.864068865776062, -0.9851729273796082),
    @(1.0256427526474, -4.58713960647583, 0.0221438650041818),
    @(2.594497442245483, 0.8452847599983215, -2.786614418029785),
    @(2.057699680328369, 1.456302762031555, -0.2620611786842346),
    @(-0.5109887719154358, -1.112385630607605, -1.43751859664917),
    @(-0.303600013256073, 2.510961532592773, -0.3182607889175415),
    @(-0.5890268087387085, -0.6258314251899719, 0.7171558141708374),
    @(1.278999090194702, 0.4355469346046448, 0.9905179738998412),
    @(-1.44057297706604, -0.845132052898407, -0.7629706859588623),
    @(-1.023504734039307, -0.6563746929168701, 0.645684540271759),
    @(-1.10917854309082, -0.5216789245605469, 0.8017606139183044),
    @(-3.361896991729736, 3.937331914901733, 2.076178312301636),
    @(2.992323398590088, 0.9390525817871094, 0.3888157308101654),
    @(-0.3591887652873993, 1.687973380088806, 0.8848382830619812),
    @(1.22615921497345, 0.5057964324951172, 0.3527746796607971),
    @(-0.1565342247486114, -0.5474879741668701, 0.3535382747650146),
    @(0.4882340431213379, -7.02800464630127, -4.413654327392578),
    @(1.305877208709717, -2.09221339225769, -1.915520668029785),
    @(1.73409366607666, -3.325092315673828, -2.228436470031738),
    @(-2.827542543411255, -0.811687171459198, -1.57450520992279)
)

$rowCount = $data.Count
for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# The old table extended one row further (to row 22); clear that leftover row.
$ws.Range("A22:C22").Clear()
